$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Function"
$ws.Range("A2").Value = "Login"
$ws.Range("B1").Value = "Note"
$ws.Range("B2").Value = "Login user with password"

$ws.Range("B3").Select()
